# Updated cryptos list on Thu Oct 10 10:13:17 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.014.94"
$ws.Range("E2").Value = "  -1.90%  "
$ws.Range("D3").Value = "2.417.22"
$ws.Range("E3").Value = "  -1.16%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "568.79"
$ws.Range("D5").NumberFormat = "General"
$ws.Range("E5").Value = "  -2.54%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "139.01"
$ws.Range("D6").NumberFormat = "General"
$ws.Range("E6").Value = "  -2.92%  "
$ws.Range("E7").Value = "  +0.15%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.527"
$ws.Range("D8").NumberFormat = "General"
$ws.Range("E8").Value = "  -1.14%  "
$ws.Range("D9").Value = "2.400.71"
$ws.Range("E9").Value = "  -1.63%  "
$ws.Range("E10").Value = "  -2.34%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.160"
$ws.Range("D11").NumberFormat = "General"
$ws.Range("E11").Value = "  -0.49%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.06"
$ws.Range("D12").NumberFormat = "General"
$ws.Range("E12").Value = "  -2.86%  "
$ws.Range("E13").Value = "  -1.82%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "25.94"
$ws.Range("D14").NumberFormat = "General"
$ws.Range("E14").Value = "  -2.06%  "
$ws.Range("B15").Value = "ShibaInu"
$ws.Range("C15").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0000171"
$ws.Range("D15").NumberFormat = "General"
$ws.Range("E15").Value = "  -2.82%  "
$ws.Range("B16").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C16").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D16").Value = "2.849.80"
$ws.Range("E16").Value = "  -0.65%  "
$ws.Range("D17").Value = "60.837.84"
$ws.Range("E17").Value = "  -1.92%  "
$ws.Range("D18").Value = "2.407.76"
$ws.Range("E18").Value = "  -1.07%  "
$ws.Range("E19").Value = "  +7.86%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "10.60"
$ws.Range("D20").NumberFormat = "General"
$ws.Range("E20").Value = "  -1.49%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "322.77"
$ws.Range("D21").NumberFormat = "General"
$ws.Range("E21").Value = "  -1.14%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.03"
$ws.Range("D22").NumberFormat = "General"
$ws.Range("E22").Value = "  -1.76%  "
$ws.Range("E23").Value = "  +1.64%  "
$ws.Range("E24").Value = "  +0.07%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.82"
$ws.Range("D25").NumberFormat = "General"
$ws.Range("E25").Value = "  -4.58%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "64.74"
$ws.Range("D26").NumberFormat = "General"
$ws.Range("E26").Value = "  -1.46%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "581.80"
$ws.Range("D27").NumberFormat = "General"
$ws.Range("E27").Value = "  -3.49%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.23"
$ws.Range("D28").NumberFormat = "General"
$ws.Range("E28").Value = "  -9.84%  "
$ws.Range("D30").Value = "0.0₃0921"
$ws.Range("E30").Value = "  -4.60%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.82"
$ws.Range("D31").NumberFormat = "General"
$ws.Range("E31").Value = "  -2.03%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.34"
$ws.Range("D32").NumberFormat = "General"
$ws.Range("E32").Value = "  -5.26%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.82"
$ws.Range("D33").NumberFormat = "General"
$ws.Range("E33").Value = "  -4.08%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.132"
$ws.Range("D34").NumberFormat = "General"
$ws.Range("E34").Value = "  -2.65%  "
$ws.Range("E35").Value = "  -0.05%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "152.04"
$ws.Range("D36").NumberFormat = "General"
$ws.Range("E36").Value = "  -0.56%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.40"
$ws.Range("D37").NumberFormat = "General"
$ws.Range("E37").Value = "  -1.74%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.59"
$ws.Range("D38").NumberFormat = "General"
$ws.Range("E38").Value = "  -6.24%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.367"
$ws.Range("D39").NumberFormat = "General"
$ws.Range("E39").Value = "  -2.51%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "18.18"
$ws.Range("D40").NumberFormat = "General"
$ws.Range("E40").Value = "  -1.32%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.13"
$ws.Range("D41").NumberFormat = "General"
$ws.Range("E41").Value = "  -2.90%  "
$ws.Range("E42").Value = "  +0.06%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.67"
$ws.Range("D43").NumberFormat = "General"
$ws.Range("E43").Value = "  -2.63%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "41.12"
$ws.Range("D44").NumberFormat = "General"
$ws.Range("E44").Value = "  -4.76%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.28"
$ws.Range("D45").NumberFormat = "General"
$ws.Range("E45").Value = "  -8.49%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "142.20"
$ws.Range("D46").NumberFormat = "General"
$ws.Range("E46").Value = "  -0.06%  "
$ws.Range("D47").Value = "0.0₆0266"
$ws.Range("E47").Value = "  +0.12%  "
$ws.Range("E48").Value = "  -3.93%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.585"
$ws.Range("D49").NumberFormat = "General"
$ws.Range("E49").Value = "  -2.53%  "
$ws.Range("B50").Value = "InjectiveProtocol"
$ws.Range("C50").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "19.34"
$ws.Range("D50").NumberFormat = "General"
$ws.Range("E50").Value = "  -2.39%  "
$ws.Range("B51").Value = "Hedera"
$ws.Range("C51").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0500"
$ws.Range("D51").NumberFormat = "General"
$ws.Range("E51").Value = "  -3.93%  "
